{"js": "// Lab 7 bug fix: correct the four p-value entries in the PCA results table.\n//   Genotype / PCA1: 0.926   -> 0.929\n//   Genotype / PCA2: 3.54e-10 -> 1.27e-10\n//   Time     / PCA1: 2.40e-31 -> 2.52e-29\n//   Time     / PCA2: 0.428   -> 0.427\nconst body = context.document.body;\n\nconst replacements = [\n  [\"0.926\", \"0.929\"],\n  [\"3.54e-10\", \"1.27e-10\"],\n  [\"2.40e-31\", \"2.52e-29\"],\n  [\"0.428\", \"0.427\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Replace only the first (and, for this document, only) match so that\n    // the run's existing character formatting (bold, font) is preserved.\n    results.items[0].insertText(replace, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lab 7 bug fix: correct the four p-value entries in the PCA results table.\n#   Genotype / PCA1: 0.926    -> 0.929\n#   Genotype / PCA2: 3.54e-10 -> 1.27e-10\n#   Time     / PCA1: 2.40e-31 -> 2.52e-29\n#   Time     / PCA2: 0.428    -> 0.427\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"0.926\", \"0.929\"),\n    @(\"3.54e-10\", \"1.27e-10\"),\n    @(\"2.40e-31\", \"2.52e-29\"),\n    @(\"0.428\", \"0.427\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceOne = 1 -- replace only the single, unique\n    # occurrence so the existing run formatting (bold) is left untouched.\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n"}
